$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 57: trade finished -> status becomes DONE and a finalized date is set
# ---------------------------------------------------------------------------
$ws.Range("H57").Value = "DONE"
$ws.Range("I57").Value2 = 42871.3746875

# ---------------------------------------------------------------------------
# Row 58 (new): LTC sell that closes the trade started on row 57
# ---------------------------------------------------------------------------
$ws.Range("A57").Copy($ws.Range("A58"))
$ws.Range("A58").Value2 = 42871.3746875

$ws.Range("B54").Copy($ws.Range("B58"))          # "Sell" (red rich text, shared string 6)

$ws.Range("C57").Copy($ws.Range("C58"))          # "        LTC" (shared string 48)

$ws.Range("D57").Copy($ws.Range("D58"))
$ws.Range("D58").Value2 = 24.5

$ws.Range("E57").Copy($ws.Range("E58"))
$ws.Range("E58").Value = "            25 USDT"

$ws.Range("F57").Copy($ws.Range("F58"))          # "        1.390 LTC" (shared string 246)

$ws.Range("G57").Copy($ws.Range("G58"))          # " LTC/USDT0000002" (shared string 249)

$ws.Range("H57").Copy($ws.Range("H58"))
$ws.Range("H58").Value = "DONE"

$ws.Range("A58").Copy($ws.Range("I58"))
$ws.Range("I58").Value2 = 42872.466319444444

$ws.Range("J54").Copy($ws.Range("J58"))
$ws.Range("J58").Value = "0.05704681 USDT (0.15%)"

$ws.Range("K54").Copy($ws.Range("K58"))
$ws.Range("K58").Value = "    ~5%"

$ws.Range("L54").Copy($ws.Range("L58"))
$ws.Range("L58").Value = " 2 day"

# ---------------------------------------------------------------------------
# Row 59 (new): XRP buy
# ---------------------------------------------------------------------------
$ws.Range("A55").Copy($ws.Range("A59"))
$ws.Range("A59").Value2 = 42873.277673611112

$ws.Range("B57").Copy($ws.Range("B59"))          # "Buy" (shared string 11 / engine collapses to matching text)

$ws.Range("C53").Copy($ws.Range("C59"))          # "        XRP" (shared string 107)

$ws.Range("D55").Copy($ws.Range("D59"))
$ws.Range("D59").Value = "              0.335
"

$ws.Range("E50").Copy($ws.Range("E59"))
$ws.Range("E59").Value = "         0.335  USDT"

$ws.Range("F50").Copy($ws.Range("F59"))
$ws.Range("F59").Value = "         111 XRP"

$ws.Range("G50").Copy($ws.Range("G59"))
$ws.Range("G59").Value = " XRP/USDT0000007"

$ws.Range("H58").Copy($ws.Range("H59"))
$ws.Range("H59").Value = "DONE"

$ws.Range("I55").Copy($ws.Range("I59"))
$ws.Range("I59").Value2 = 42873.277673611112

$ws.Range("K50").Copy($ws.Range("K59"))          # "     " (shared string 171)

# ---------------------------------------------------------------------------
# Row 60 (new): XRP sell, still in progress
# ---------------------------------------------------------------------------
$ws.Range("A59").Copy($ws.Range("A60"))
$ws.Range("A60").Value2 = 42873.277673611112

$ws.Range("B54").Copy($ws.Range("B60"))          # "Sell" (red rich text, shared string 6)

$ws.Range("C59").Copy($ws.Range("C60"))          # "        XRP" (shared string 107)

$ws.Range("D56").Copy($ws.Range("D60"))
$ws.Range("D60").Value = "              0.333
"

$ws.Range("E59").Copy($ws.Range("E60"))
$ws.Range("E60").Value = "         0.365  USDT"

$ws.Range("F59").Copy($ws.Range("F60"))
$ws.Range("F60").Value = "         111 XRP"

$ws.Range("G59").Copy($ws.Range("G60"))
$ws.Range("G60").Value = " XRP/USDT0000007"

$ws.Range("H57").Copy($ws.Range("H60"))
$ws.Range("H60").Value = "IN PROGRESS"

$ws.Range("I59").Copy($ws.Range("I60"))
$ws.Range("I60").ClearContents()

$ws.Range("K59").Copy($ws.Range("K60"))          # "     " (shared string 171)

# ---------------------------------------------------------------------------
# Selection moved to F65
# ---------------------------------------------------------------------------
$ws.Range("F65").Select()
